$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text in the source data
# (e.g. "569.22", "62.948.39"). Excel auto-detects numeric-looking strings
# and would silently convert them to numbers (losing formatting / exact
# text such as trailing zeros or thousand-dot formatting). Forcing the
# column to Text format first preserves every Price cell as a string,
# matching the original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.898.87"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "2.433.67"
$ws.Range("E3").Value = "  -0.37%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "569.03"
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("D6").Value = "145.93"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "0.536"
$ws.Range("E8").Value = "  +0.63%  "

$ws.Range("D9").Value = "2.428.41"
$ws.Range("E9").Value = "  -0.68%  "

$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("E11").Value = "  +1.08%  "

$ws.Range("D12").Value = "5.22"
$ws.Range("E12").Value = "  -1.68%  "

$ws.Range("D13").Value = "0.351"
$ws.Range("E13").Value = "  -0.65%  "

$ws.Range("D14").Value = "26.90"
$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("D15").Value = "0.0000178"
$ws.Range("E15").Value = "  -2.53%  "

$ws.Range("D16").Value = "2.871.76"
$ws.Range("E16").Value = "  -0.42%  "

$ws.Range("D17").Value = "63.314.03"
$ws.Range("E17").Value = "  +1.23%  "

$ws.Range("D18").Value = "2.419.07"
$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("D19").Value = "11.24"
$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("D20").Value = "7.30"
$ws.Range("E20").Value = "  +5.24%  "

$ws.Range("D21").Value = "325.30"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").Value = "4.16"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("D23").Value = "2.07"
$ws.Range("E23").Value = "  +12.28%  "

$ws.Range("D24").Value = "0.997"
$ws.Range("E24").Value = "  -0.31%  "

$ws.Range("D25").Value = "65.05"
$ws.Range("E25").Value = "  -3.38%  "

$ws.Range("D26").Value = "614.57"
$ws.Range("E26").Value = "  +6.42%  "

$ws.Range("D27").Value = "8.84"
$ws.Range("E27").Value = "  +1.11%  "

$ws.Range("D28").Value = "0.0000102"
$ws.Range("E28").Value = "  +0.42%  "

$ws.Range("D29").Value = "2.578.49"
$ws.Range("E29").Value = "  +0.63%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "1.49"
$ws.Range("E31").Value = "  +2.92%  "

$ws.Range("D32").Value = "8.21"
$ws.Range("E32").Value = "  -2.74%  "

$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "0.141"
$ws.Range("E33").Value = "  -4.60%  "

$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "1.89"
$ws.Range("E34").Value = "  +0.50%  "

$ws.Range("D35").Value = "5.15"
$ws.Range("E35").Value = "  +6.25%  "

$ws.Range("D36").Value = "1.51"
$ws.Range("E36").Value = "  -2.77%  "

$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("E38").Value = "  -1.51%  "

$ws.Range("D39").Value = "5.40"
$ws.Range("E39").Value = "  -0.47%  "

$ws.Range("D40").Value = "18.64"
$ws.Range("E40").Value = "  -0.91%  "

$ws.Range("D41").Value = "145.63"
$ws.Range("E41").Value = "  -1.73%  "

$ws.Range("D42").Value = "2.66"
$ws.Range("E42").Value = "  +9.53%  "

$ws.Range("E43").Value = "  -2.29%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").Value = "147.55"
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("D46").Value = "3.72"
$ws.Range("E46").Value = "  +1.48%  "

$ws.Range("D47").Value = "21.09"
$ws.Range("E47").Value = "  +2.62%  "

$ws.Range("D48").Value = "0.0533"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("D49").Value = "0.595"
$ws.Range("E49").Value = "  -1.20%  "

$ws.Range("D50").Value = "0.0233"
$ws.Range("E50").Value = "  +0.72%  "

$ws.Range("D51").Value = "0.0911"
$ws.Range("E51").Value = "  -1.48%  "
